$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that receive new values for rows 2-25: B, C, D, F, G, L, M
$cols = @("B", "C", "D", "F", "G", "L", "M")

# New values per row (row index 2..25), in the order of $cols
$data = @(
    @(1.746658609027293, 0.1294864946956409, 0.2382925647604992, 3.97516721055527, 0.002548948803578662, 0.2299822569156404, 0.3615611641911798),
    @(1.679612626566666, 0.112865639352151, 0.2278923665241877, 3.773134890883966, 0.002555607924258156, 0.2273598256816172, 0.3514363342560571),
    @(1.640167408767638, 0.1026582531967222, 0.2215020269047301, 3.650072339435667, 0.0025599050354407, 0.2258424243402573, 0.3455206123746919),
    @(1.624523524921983, 0.09849726812180393, 0.2188962178649518, 3.600163069654712, 0.002561708748109839, 0.2252473342195955, 0.3431852579160761),
    @(1.621951791925255, 0.09780623071162609, 0.2184634071517166, 3.591889912510453, 0.002562011436447435, 0.2251499230704681, 0.34280201561333),
    @(1.639954690559364, 0.1026021434634288, 0.2214668916270028, 3.64939828579395, 0.002559929147765878, 0.225834304641765, 0.345488812304005),
    @(1.723182368358437, 0.1237555363234151, 0.2347071957291291, 3.905296753111031, 0.002551201741956797, 0.2290587432583351, 0.3580074518877154),
    @(1.900175264635152, 0.1652615572905063, 0.2606606970623204, 4.415336209687723, 0.002535731378747387, 0.2361217084207397, 0.3849616878855784),
    @(2.038810719090975, 0.1958309218951797, 0.2797592966107061, 4.795667087927256, 0.002525354380312017, 0.2417678798542084, 0.4062584581069615),
    @(2.103791561942558, 0.2097673853499771, 0.2884625404485064, 4.970038155955933, 0.002520845553743459, 0.244437087620156, 0.4162774715641504),
    @(2.128676894214152, 0.2150501809357479, 0.2917610551484984, 5.036273146222129, 0.00251916840491568, 0.2454624330579378, 0.4201194805291735),
    @(2.123304955372191, 0.2139121818594276, 0.2910505274902846, 5.021999012144704, 0.002519528266863119, 0.245240956931724, 0.4192898922258337),
    @(2.105833293375554, 0.2102018909013168, 0.2887338519552713, 4.975483196757807, 0.002520706968295537, 0.2445211507963165, 0.416592590941633),
    @(2.095167759330707, 0.2079299585816443, 0.287315201942107, 4.947017822930661, 0.002521432890444342, 0.244082149449099, 0.4149466833586928),
    @(2.034602866688431, 0.1949208361917556, 0.2791908760582942, 4.784299624968014, 0.002525653286202467, 0.24159547543222, 0.405610388951068),
    @(1.997940917201788, 0.186948628995367, 0.2742111818458284, 4.684831367110263, 0.002528296447565823, 0.2400958445970218, 0.3999679532822995),
    @(1.977033977448116, 0.1823660176818009, 0.2713483873641849, 4.627746980469794, 0.002529836664295679, 0.239242772888943, 0.3967537029443378),
    @(1.969986075460724, 0.1808148734927499, 0.2703793155868937, 4.60844073880196, 0.002530361586096088, 0.2389555623720838, 0.395670746704667),
    @(2.001824988013595, 0.1877969879738259, 0.2747411302056264, 4.695406708698698, 0.002528013016222627, 0.2402545014666657, 0.400565375360749),
    @(2.110957560180225, 0.2112915397596566, 0.2894142353112272, 4.989140395132836, 0.002520359935836445, 0.2447321790413639, 0.4173835468219167),
    @(2.183906826936095, 0.2266783057142447, 0.2990204754385388, 5.182308495133327, 0.002515534411434364, 0.2477435940194113, 0.4286553349214515),
    @(2.144822696304232, 0.2184628605625392, 0.2938917327364834, 5.079098461045078, 0.002518093827445235, 0.2461285382017024, 0.4226135915999265),
    @(2.000068468755728, 0.1874134427220611, 0.2745015402918796, 4.690625282133595, 0.002528141091449729, 0.2401827444014515, 0.4002951884835397),
    @(1.85079843436165, 0.154024337096331, 0.2536369215640235, 4.276415986309161, 0.002539741878002749, 0.2341311744863933, 0.3774095917151783)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $cellRef = "$($cols[$j])$row"
        $ws.Range($cellRef).Value = $rowValues[$j]
    }
}
